$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "30.447.37"
$ws.Cells.Item(2, 5).Value = "  +0.53%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.866.96"
$ws.Cells.Item(3, 5).Value = "  -0.05%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.07%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "235.62"
$ws.Cells.Item(5, 5).Value = "  -0.87%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(6, 5).Value = "  -0.08%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4786"
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2793"
$ws.Cells.Item(8, 5).Value = "  -0.36%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06537"
$ws.Cells.Item(9, 5).Value = "  +0.74%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.875.61"
$ws.Cells.Item(10, 5).Value = "  +0.52%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07444"
$ws.Cells.Item(11, 5).Value = "  +0.09%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "16.22"
$ws.Cells.Item(12, 5).Value = "  -1.64%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.071"
$ws.Cells.Item(13, 5).Value = "  +0.17%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "86.93"
$ws.Cells.Item(14, 5).Value = "  -1.18%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.6395"
$ws.Cells.Item(15, 5).Value = "  -2.08%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "30.438.19"
$ws.Cells.Item(16, 5).Value = "  +0.57%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.9999"
$ws.Cells.Item(17, 5).Value = "  -0.09%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -2.11%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "233.29"
$ws.Cells.Item(19, 5).Value = "  +6.61%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.000007465"
$ws.Cells.Item(20, 5).Value = "  -1.22%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "BinanceUSD"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.000"
$ws.Cells.Item(21, 5).Value = "  -0.12%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.123"
$ws.Cells.Item(22, 5).Value = "  -2.72%  "

# Row 23
$ws.Cells.Item(23, 2).Value = "Chainlink"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.097"
$ws.Cells.Item(23, 5).Value = "  -0.79%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "Monero"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "168.37"
$ws.Cells.Item(24, 5).Value = "  -0.16%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "Cosmos"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.306"
$ws.Cells.Item(25, 5).Value = "  +0.07%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "18.15"
$ws.Cells.Item(26, 5).Value = "  -1.47%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "LidoDAOToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.897"
$ws.Cells.Item(27, 5).Value = "  -3.55%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "Stellar"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.1048"
$ws.Cells.Item(28, 5).Value = "  +11.90%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.379"
$ws.Cells.Item(29, 5).Value = "  -5.14%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.272"
$ws.Cells.Item(30, 5).Value = "  -0.56%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.968"
$ws.Cells.Item(31, 5).Value = "  -0.99%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Hedera"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.04980"
$ws.Cells.Item(32, 5).Value = "  -1.52%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "ARBITRUM"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.171"
$ws.Cells.Item(33, 5).Value = "  -2.51%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "ImmutableX"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.7390"
$ws.Cells.Item(34, 5).Value = "  -1.08%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Frax"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9992"
$ws.Cells.Item(35, 5).Value = "  +0.07%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "HuobiToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.708"
$ws.Cells.Item(36, 5).Value = "  +0.01%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.01946"
$ws.Cells.Item(37, 5).Value = "  +7.34%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "MXToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.637"
$ws.Cells.Item(38, 5).Value = "  +0.95%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "TrustWalletToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.9123"
$ws.Cells.Item(39, 5).Value = "  +0.87%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "RenderToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.040"
$ws.Cells.Item(40, 5).Value = "  -1.71%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "Quant"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "106.32"
$ws.Cells.Item(41, 5).Value = "  -0.42%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "PaxDollar"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.9964"
$ws.Cells.Item(42, 5).Value = "  -0.64%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "TheSandbox"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.4175"
$ws.Cells.Item(43, 5).Value = "  -2.00%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.564"
$ws.Cells.Item(44, 5).Value = "  -6.49%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Aptos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "7.159"
$ws.Cells.Item(45, 5).Value = "  -2.86%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "61.77"
$ws.Cells.Item(46, 5).Value = "  -3.14%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Algorand"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.1223"
$ws.Cells.Item(47, 5).Value = "  -4.58%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "8.857"
$ws.Cells.Item(48, 5).Value = "  -1.34%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Elrond"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "33.49"
$ws.Cells.Item(49, 5).Value = "  -0.33%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -4.03%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.05638"
$ws.Cells.Item(51, 5).Value = "  +0.33%  "
